$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the font weight on rows 3-5 so they pick up the plain
# (non-fill) cell styles instead of the legacy fill-flagged ones.
$ws.Range("A3:A5").Font.Bold = $true
$ws.Range("B3:B5").Font.Bold = $false

# Insert a new row above the "Ave flux um/m2" row and record the
# average air pressure measurement there.
$ws.Rows(33).Insert()
$ws.Range("A33").Value = "Ave Pressure kpa (air)"
$ws.Range("B33").Value = 62.4

# Move the active selection to the newly added cell.
$ws.Range("B33").Select() | Out-Null
